# Mise à jour de l'application
# Append 18 new training-log rows (874-891) to Feuil1, mirroring the
# existing row layout/format (copy formats down from the last existing
# row, then fill in values + the I = C*D "Charge" formula).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

$lastRow = 873
$firstNew = 874
$lastNew = 891

# Clone the formatting (styles / number formats / column widths behaviour)
# of the last existing data row down across the new block so the new
# cells pick up the same styles (date format in A, name font in B-H, etc.)
$ws.Range("A${lastRow}:I${lastRow}").Copy() | Out-Null
$ws.Range("A${firstNew}:I${lastNew}").PasteSpecial(-4122) | Out-Null

# New rows recorded on 2026-02-23 (serial date 46076)
$rows = @(
    @{ Row=874; Name="Yoann Martelat";   C=60; D=5; E=6; F=3; G="Genou";             H=7  },
    @{ Row=875; Name="Kamal Bafounta";   C=60; D=5; E=2; F=0; G="";                  H=5  },
    @{ Row=876; Name="Omar Benyounes";   C=60; D=4; E=6; F=0; G="";                  H=0  },
    @{ Row=877; Name="Ryad Kralladi";    C=60; D=5; E=4; F=0; G="";                  H=3  },
    @{ Row=878; Name="Naim Ighbane";     C=60; D=8; E=3; F=3; G="";                  H=3  },
    @{ Row=879; Name="Mehdi Boussaid";   C=60; D=6; E=6; F=0; G="";                  H=7  },
    @{ Row=880; Name="Sofiane Belle";    C=60; D=5; E=4; F=1; G="Ventre";            H=5  },
    @{ Row=881; Name="Maé Clavel";       C=60; D=4; E=3; F=2; G="Tibia";             H=4  },
    @{ Row=882; Name="Theo Owono";       C=60; D=5; E=5; F=0; G="";                  H=7  },
    @{ Row=883; Name="Nathanael Beta";   C=60; D=6; E=8; F=0; G="";                  H=8  },
    @{ Row=884; Name="Levy Ndoutoume";   C=60; D=5; E=5; F=0; G="";                  H=7  },
    @{ Row=885; Name="Ilan Ihaddadene";  C=60; D=8; E=6; F=0; G="";                  H=10 },
    @{ Row=886; Name="Mattheo Haon";     C=60; D=8; E=0; F=0; G="";                  H=8  },
    @{ Row=887; Name="Karahali Souaré";  C=60; D=6; E=6; F=6; G="";                  H=1  },
    @{ Row=888; Name="Romain Thunet";    C=60; D=6; E=4; F=3; G="Mollet béquille ";  H=6  },
    @{ Row=889; Name="Naim Dhib";        C=60; D=5; E=2; F=2; G="Psoas";             H=3  },
    @{ Row=890; Name="Yoan Zouma";       C=60; D=7; E=3; F=3; G="Cheville";          H=7  },
    @{ Row=891; Name="Jeremie Laurent";  C=60; D=6; E=5; F=0; G="";                  H=7  }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Range("A$row").Value = 46076
    $ws.Range("B$row").Value = $r.Name
    $ws.Range("C$row").Value = $r.C
    $ws.Range("D$row").Value = $r.D
    $ws.Range("E$row").Value = $r.E
    $ws.Range("F$row").Value = $r.F
    if ($r.G -ne "") {
        $ws.Range("G$row").Value = $r.G
    }
    $ws.Range("H$row").Value = $r.H
}

# Charge = Volume * Intensité, same formula used throughout column I
$ws.Range("I${firstNew}:I${lastNew}").Formula = "=C${firstNew}*D${firstNew}"

# Reflect the updated view position/selection from the edited workbook
$ws.Application.Goto($ws.Range("A865"), $false)
$ws.Range("C894").Select() | Out-Null
